# The deck's slide-master theme ("Integral" colours) is switched over to
# the default Office palette, matching the author's re-colour of the
# design's theme colour scheme (Design > Variants > Colors > "Office").
#
# PowerPoint keeps the twelve core theme colours on the master's
# ColorScheme collection (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink,
# in that fixed index order) - reassigning each one rewrites the
# underlying theme part's <a:clrScheme> with the new palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

function New-ComRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette: the standard Office theme colours.
$scheme.Colors(1).RGB  = New-ComRgb 0x00 0x00 0x00   # dk1
$scheme.Colors(2).RGB  = New-ComRgb 0xFF 0xFF 0xFF   # lt1
$scheme.Colors(3).RGB  = New-ComRgb 0x44 0x54 0x6A   # dk2
$scheme.Colors(4).RGB  = New-ComRgb 0xE7 0xE6 0xE6   # lt2
$scheme.Colors(5).RGB  = New-ComRgb 0x5B 0x9B 0xD5   # accent1
$scheme.Colors(6).RGB  = New-ComRgb 0xED 0x7D 0x31   # accent2
$scheme.Colors(7).RGB  = New-ComRgb 0xA5 0xA5 0xA5   # accent3
$scheme.Colors(8).RGB  = New-ComRgb 0xFF 0xC0 0x00   # accent4
$scheme.Colors(9).RGB  = New-ComRgb 0x44 0x72 0xC4   # accent5
$scheme.Colors(10).RGB = New-ComRgb 0x70 0xAD 0x47   # accent6
$scheme.Colors(11).RGB = New-ComRgb 0x05 0x63 0xC1   # hlink
$scheme.Colors(12).RGB = New-ComRgb 0x95 0x4F 0x72   # folHlink
